$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC!row116
$ws_ALC.Range("H116").Value = 2218.3157
$ws_ALC.Range("I116").Value = 1999.091
$ws_ALC.Range("J116").Value = 2519.75
$ws_ALC.Range("K116").Value = 1999.091
$ws_ALC.Range("L116").Value = 2519.75
$ws_ALC.Range("M116").Value = 1442.909
$ws_ALC.Range("N116").Value = -9403.75

# ALC!row138
$ws_ALC.Range("H138").Value = 3587.375
$ws_ALC.Range("I138").Value = 2197.2307
$ws_ALC.Range("J138").Value = 3941.7256
$ws_ALC.Range("K138").Value = 6591.6921
$ws_ALC.Range("L138").Value = 11825.1768
$ws_ALC.Range("M138").Value = -1451.6921
$ws_ALC.Range("N138").Value = -22105.1768

# ALC!row141
$ws_ALC.Range("H141").Value = 2803.8
$ws_ALC.Range("I141").Value = 1717.1111
$ws_ALC.Range("J141").Value = 5060.769
$ws_ALC.Range("K141").Value = 5151.3333
$ws_ALC.Range("L141").Value = 15182.307
$ws_ALC.Range("M141").Value = 28.66669999999976
$ws_ALC.Range("N141").Value = -25542.307

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM!row32
$ws_ARM.Range("H32").Value = 17466.123
$ws_ARM.Range("I32").Value = 18833.406
$ws_ARM.Range("K32").Value = 18833.406
$ws_ARM.Range("M32").Value = -18546.406

# ARM!row35
$ws_ARM.Range("H35").Value = 2248
$ws_ARM.Range("I35").Value = 2248
$ws_ARM.Range("K35").Value = 2248
$ws_ARM.Range("M35").Value = -1842

# ARM!row45
$ws_ARM.Range("H45").Value = 2190.4517
$ws_ARM.Range("I45").Value = 1954.3334
$ws_ARM.Range("J45").Value = 3000
$ws_ARM.Range("K45").Value = 1954.3334
$ws_ARM.Range("L45").Value = 3000
$ws_ARM.Range("M45").Value = -1577.3334
$ws_ARM.Range("N45").Value = -3754

# ARM!row61
$ws_ARM.Range("H61").Value = 8062.2354
$ws_ARM.Range("I61").Value = 4944.826
$ws_ARM.Range("J61").Value = 14580.454
$ws_ARM.Range("K61").Value = 4944.826
$ws_ARM.Range("L61").Value = 14580.454
$ws_ARM.Range("M61").Value = -4732.826
$ws_ARM.Range("N61").Value = -15004.454

# ARM!row136
$ws_ARM.Range("H136").Value = 8062.2354
$ws_ARM.Range("I136").Value = 4944.826
$ws_ARM.Range("J136").Value = 14580.454
$ws_ARM.Range("K136").Value = 14834.478
$ws_ARM.Range("L136").Value = 43741.362
$ws_ARM.Range("M136").Value = -12284.478
$ws_ARM.Range("N136").Value = -48841.362

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM!row5
$ws_BSM.Range("H5").Value = 3399.5
$ws_BSM.Range("J5").Value = 9502.5
$ws_BSM.Range("L5").Value = 9502.5
$ws_BSM.Range("N5").Value = -9728.5

# BSM!row19
$ws_BSM.Range("H19").Value = 15000
$ws_BSM.Range("J19").Value = 15000
$ws_BSM.Range("L19").Value = 15000
$ws_BSM.Range("N19").Value = -15346

# BSM!row134
$ws_BSM.Range("H134").Value = 44884.39
$ws_BSM.Range("I134").Value = 1515.65
$ws_BSM.Range("K134").Value = 4546.950000000001
$ws_BSM.Range("M134").Value = -2011.950000000001

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP!row7
$ws_CRP.Range("H7").Value = 82.71429000000001
$ws_CRP.Range("I7").Value = 76
$ws_CRP.Range("J7").Value = 99.5
$ws_CRP.Range("K7").Value = 76
$ws_CRP.Range("L7").Value = 99.5
$ws_CRP.Range("M7").Value = 37
$ws_CRP.Range("N7").Value = -325.5

# CRP!row31
$ws_CRP.Range("H31").Value = 811425.3
$ws_CRP.Range("I31").Value = 12072.889
$ws_CRP.Range("J31").Value = 1325294.8
$ws_CRP.Range("K31").Value = 12072.889
$ws_CRP.Range("L31").Value = 1325294.8
$ws_CRP.Range("M31").Value = -11777.889
$ws_CRP.Range("N31").Value = -1325884.8

# CRP!row34
$ws_CRP.Range("H34").Value = 811425.3
$ws_CRP.Range("I34").Value = 12072.889
$ws_CRP.Range("J34").Value = 1325294.8
$ws_CRP.Range("K34").Value = 12072.889
$ws_CRP.Range("L34").Value = 1325294.8
$ws_CRP.Range("M34").Value = -11870.889
$ws_CRP.Range("N34").Value = -1325698.8

# CRP!row132
$ws_CRP.Range("H132").Value = 3062.0278
$ws_CRP.Range("I132").Value = 2722.1785
$ws_CRP.Range("J132").Value = 4251.5
$ws_CRP.Range("K132").Value = 8166.5355
$ws_CRP.Range("L132").Value = 12754.5
$ws_CRP.Range("M132").Value = -5636.5355
$ws_CRP.Range("N132").Value = -17814.5

# CRP!row135
$ws_CRP.Range("H135").Value = 54247.855
$ws_CRP.Range("J135").Value = 54247.855
$ws_CRP.Range("L135").Value = 54247.855
$ws_CRP.Range("N135").Value = -64387.855

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL!row3
$ws_CUL.Range("H3").Value = 3539.6562
$ws_CUL.Range("I3").Value = 2260
$ws_CUL.Range("J3").Value = 4989.933
$ws_CUL.Range("K3").Value = 6780
$ws_CUL.Range("L3").Value = 14969.799
$ws_CUL.Range("M3").Value = -6668
$ws_CUL.Range("N3").Value = -15193.799

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM!row97
$ws_GSM.Range("H97").Value = 1362
$ws_GSM.Range("I97").Value = 1510.2727
$ws_GSM.Range("J97").Value = 1035.8
$ws_GSM.Range("K97").Value = 1510.2727
$ws_GSM.Range("L97").Value = 1035.8
$ws_GSM.Range("M97").Value = -1014.2727
$ws_GSM.Range("N97").Value = -2027.8

# GSM!row132
$ws_GSM.Range("H132").Value = 5295.7666
$ws_GSM.Range("I132").Value = 1343.05
$ws_GSM.Range("J132").Value = 13201.2
$ws_GSM.Range("K132").Value = 4029.15
$ws_GSM.Range("L132").Value = 39603.60000000001
$ws_GSM.Range("M132").Value = -1499.15
$ws_GSM.Range("N132").Value = -44663.60000000001

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW!row22
$ws_LTW.Range("H22").Value = 449.75
$ws_LTW.Range("I22").Value = 449.75
$ws_LTW.Range("J22").Value = 0
$ws_LTW.Range("K22").Value = 449.75
$ws_LTW.Range("L22").Value = 0
$ws_LTW.Range("M22").Value = -154.75
$ws_LTW.Range("N22").ClearContents()

# LTW!row27
$ws_LTW.Range("H27").Value = 449.75
$ws_LTW.Range("I27").Value = 449.75
$ws_LTW.Range("J27").Value = 0
$ws_LTW.Range("K27").Value = 449.75
$ws_LTW.Range("L27").Value = 0
$ws_LTW.Range("M27").Value = -342.75
$ws_LTW.Range("N27").ClearContents()

# LTW!row132
$ws_LTW.Range("H132").Value = 5385.0454
$ws_LTW.Range("I132").Value = 5925.967
$ws_LTW.Range("J132").Value = 4225.9287
$ws_LTW.Range("K132").Value = 17777.901
$ws_LTW.Range("L132").Value = 12677.7861
$ws_LTW.Range("M132").Value = -15247.901
$ws_LTW.Range("N132").Value = -17737.7861

# LTW!row136
$ws_LTW.Range("H136").Value = 2399.388
$ws_LTW.Range("I136").Value = 1613.1
$ws_LTW.Range("J136").Value = 3564.2593
$ws_LTW.Range("K136").Value = 4839.299999999999
$ws_LTW.Range("L136").Value = 10692.7779
$ws_LTW.Range("M136").Value = -2289.299999999999
$ws_LTW.Range("N136").Value = -15792.7779

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR!row46
$ws_WVR.Range("H46").Value = 41357.25
$ws_WVR.Range("J46").Value = 41357.25
$ws_WVR.Range("L46").Value = 41357.25
$ws_WVR.Range("N46").Value = -41819.25

# WVR!row132
$ws_WVR.Range("H132").Value = 1691.65
$ws_WVR.Range("I132").Value = 1616.9688
$ws_WVR.Range("J132").Value = 1990.375
$ws_WVR.Range("K132").Value = 4850.9064
$ws_WVR.Range("L132").Value = 5971.125
$ws_WVR.Range("M132").Value = -2320.9064
$ws_WVR.Range("N132").Value = -11031.125

# WVR!row134
$ws_WVR.Range("H134").Value = 41357.25
$ws_WVR.Range("J134").Value = 41357.25
$ws_WVR.Range("L134").Value = 124071.75
$ws_WVR.Range("N134").Value = -129141.75

# WVR!row136
$ws_WVR.Range("H136").Value = 5968.1113
$ws_WVR.Range("I136").Value = 5370.143
$ws_WVR.Range("J136").Value = 6715.5713
$ws_WVR.Range("K136").Value = 16110.429
$ws_WVR.Range("L136").Value = 20146.7139
$ws_WVR.Range("M136").Value = -13560.429
$ws_WVR.Range("N136").Value = -25246.7139

Write-Output "Applied all Pandaemonium_Profits updates"